# Apply cryptos list update (price/volume refresh + a few row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.405.72'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '1.916.20'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.99'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4696'
$ws.Range('E7').Value = '  -2.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2845'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06832'
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '106.72'
$ws.Range('E10').Value = '  +11.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.02'
$ws.Range('E11').Value = '  -4.35%  '
$ws.Range('D12').Value = '1.909.66'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07637'
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.179'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6538'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '286.34'
$ws.Range('E16').Value = '  -4.76%  '
$ws.Range('D17').Value = '30.397.68'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007588'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.95'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').Value = '2.158.82'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.208'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.169'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.04'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.224'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.25'
$ws.Range('E27').Value = '  +6.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.029'
$ws.Range('E28').Value = '  +3.33%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.368'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.139'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.939'
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05036'
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7352'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').Value = '  -3.38%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9991'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.721'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02003'
$ws.Range('E38').Value = '  +2.48%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.672'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.043'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '108.88'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8732'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.835'
$ws.Range('E43').Value = '  +4.00%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '52.72'
$ws.Range('E45').Value = '  +25.09%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4190'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '67.22'
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.101'
$ws.Range('E48').Value = '  -4.07%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.207'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1202'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.55'
$ws.Range('E51').Value = '  -0.62%  '
